# Update test case names ("wrong" -> "invalid") across the three sheets,
# plus refresh saved selections/zoom, matching the authored diff.

$wb = $excel.ActiveWorkbook

# ---------- Login sheet ----------
$login = $wb.Worksheets.Item("Login")
$login.Range("A2").Value = "negative case - invalid phone number and password"
$login.Range("A3").Value = "negative case - invalid phone number"
$login.Range("A4").Value = "negative case - invalid password"
$login.Range("A6").Value = "negative case - invalid phone number format (not numerical)"
$login.Range("A7").Value = "negative case - invalid phone number format (exceed 12 characters)"
$login.Range("A8").Value = "negative case - invalid phone number (outside Indonesia)"
$login.Range("A9").Value = "negative case - invalid phone number (unregistered number)"
$login.Range("A10").Value = "negative case - invalid phone number (unverified number)"
$login.Range("A11").Value = "negative case - invalid password (no lower case alphabetical character)"
$login.Range("A12").Value = "negative case - invalid password (no upper case alphabetical character)"
$login.Range("A13").Value = "negative case - invalid password (no numerical character)"
$login.Range("A14").Value = "negative case - invalid password (below 8 characters)"
$login.Range("A15").Value = "negative case - invalid password (exceed 16 characters)"

# ---------- Register sheet ----------
$register = $wb.Worksheets.Item("Register")
$register.Range("A2").Value = "negative case - invalid full name, email, phone number, password, confirm password"
$register.Range("A3").Value = "negative case - invalid full name, email, phone number, password"
$register.Range("A4").Value = "negative case - invalid full name, email, phone number, confirm password"
$register.Range("A5").Value = "negative case - invalid full name, email, phone number"
$register.Range("A6").Value = "negative case - invalid full name, email, password, confirm password"
$register.Range("A7").Value = "negative case - invalid full name, email, password"
$register.Range("A8").Value = "negative case - invalid full name, email, confirm password"
$register.Range("A9").Value = "negative case - invalid full name, email"
$register.Range("A10").Value = "negative case - invalid full name, phone number, password, confirm password"
$register.Range("A11").Value = "negative case - invalid full name, phone number, password"
$register.Range("A12").Value = "negative case - invalid full name, phone number , confirm password"
$register.Range("A13").Value = "negative case - invalid full name, phone number"
$register.Range("A14").Value = "negative case - invalid full name, password, confirm password"
$register.Range("A15").Value = "negative case - invalid full name, password"
$register.Range("A16").Value = "negative case - invalid full name, confirm password"
$register.Range("A17").Value = "negative case - invalid full name"
$register.Range("A18").Value = "negative case - invalid email, phone number, password, confirm password"
$register.Range("A19").Value = "negative case - invalid email, phone number, password"
$register.Range("A20").Value = "negative case - invalid email, phone number, confirm password"
$register.Range("A21").Value = "negative case - invalid email, phone number"
$register.Range("A22").Value = "negative case - invalid email, password, confirm password"
$register.Range("A23").Value = "negative case - invalid email, password"
$register.Range("A24").Value = "negative case - invalid email, confirm password"
$register.Range("A25").Value = "negative case - invalid email"
$register.Range("A26").Value = "negative case - inalid phone number, password, confirm password"
$register.Range("A27").Value = "negative case - invalid phone number, password"
$register.Range("A28").Value = "negative case - invalid phone number, confirm password"
$register.Range("A29").Value = "negative case - invalid phone number"
$register.Range("A30").Value = "negative case - invalid password, confirm password"
$register.Range("A31").Value = "negative case - inavalid password"
$register.Range("A32").Value = "negative case - invalid confirm password"
$register.Range("A44").Value = "negative case - invalid full name (below 3 characters)"
$register.Range("A45").Value = "negative case - invalid full name (exceed 20 characters)"
$register.Range("A46").Value = "negative case - invalid full name (non alphabetical)"
$register.Range("A47").Value = "negative case - invalid email (no @)"
$register.Range("A48").Value = "negative case - invalid phone number format (not numerical)"
$register.Range("A49").Value = "negative case - invalid phone number format (exceed 12 characters)"
$register.Range("A50").Value = "negative case - invalid phone number (outside Indonesia)"
$register.Range("A51").Value = "negative case - invalid phone number (unregistered number)"
$register.Range("A52").Value = "negative case - invalid phone number (unverified number)"
$register.Range("A53").Value = "negative case - invalid password (no lower case alphabetical character)"
$register.Range("A54").Value = "negative case - invalid password (no upper case alphabetical character)"
$register.Range("A55").Value = "negative case - invalid password (no numerical character)"
$register.Range("A56").Value = "negative case - invalid password (below 8 characters)"
$register.Range("A57").Value = "negative case - invalid password (exceed 16 characters)"
$register.Range("A58").Value = "negative case - invalid confirm password (no lower case alphabetical character)"
$register.Range("A59").Value = "negative case - invalid confirm password (no upper case alphabetical character)"
$register.Range("A60").Value = "negative case - invalid confirm password (no numerical character)"
$register.Range("A61").Value = "negative case - invalid confirm password (below 8 characters)"
$register.Range("A62").Value = "negative case - invalid confirm password (exceed 16 characters)"

# ---------- Forgot Password sheet ----------
$forgot = $wb.Worksheets.Item("Forgot Password")
$forgot.Range("A2").Value = "negative case - invalid new password and email"
$forgot.Range("A3").Value = "negative case - invalid email"
$forgot.Range("A6").Value = "negative case - invalid password (no lower case alphabetical character)"
$forgot.Range("A7").Value = "negative case - invalid password (no upper case alphabetical character)"
$forgot.Range("A8").Value = "negative case - invalid password (no numerical character)"
$forgot.Range("A9").Value = "negative case - invalid password (below 8 characters)"
$forgot.Range("A10").Value = "negative case - invalid password (exceed 16 characters)"
$forgot.Range("A11").Value = "negative case - invalid email (no @)"

# ---------- Restore selections / view state ----------
$login.Activate()
$login.Range("B17").Select()

$register.Activate()
$excel.ActiveWindow.Zoom = 101
$register.Range("A30").Select()

$forgot.Activate()
$forgot.Range("A20").Select()

$forgot.Activate()
